$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared-string values used for Motor / Chasis / Patente columns (rows 2-3)
$ws.Range("H2").Value = "ZZZ520"
$ws.Range("I2").Value = "ABC0987AX316"
$ws.Range("J2").Value = "MMAA09XFGS311"
$ws.Range("H3").Value = "ZZZ521"
$ws.Range("I3").Value = "ABC0987AX317"
$ws.Range("J3").Value = "MMAA09XFGS312"

# Remove the empty styled cells in H4:J5 (no longer part of the used set)
$ws.Range("H4:J5").Clear()

# Update selection to match new active range
$ws.Range("H2:J3").Select()

# Update the conditional formatting range so it covers H6:J18 and H2:J3
# (keep the existing rule's dxf/priority instead of recreating it)
$fc = $ws.Range("H2:J18").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("H6:J18,H2:J3"))
